$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.923.67'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '2.635.03'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.549'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").Value = '2.633.12'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("E10").Value = '  +10.59%  '
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("E15").Value = '  +3.82%  '
$ws.Range("D16").Value = '3.110.17'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = '67.816.24'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '2.654.34'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '374.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.85%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("D28").Value = '2.760.29'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '575.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.19'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.16'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.41%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("E44").Value = '  +11.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("E50").Value = '  +7.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.58%  '

Write-Host "Applied cryptos update"
